# Updated cryptos list on Fri Dec 22 14:47:04 UTC 2023 with GitHub Actions
#
# This script refreshes the "Price" (D) and "Volume(1h)" (E) columns of the
# cryptos worksheet with a new scrape snapshot, and reflects two small
# ranking swaps (VeChain/Kaspa around rows 36-37, and the
# FraxShare/MultiversX/THORChain rotation around rows 44-46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text price strings (e.g. "43.603.77", "0.108") that
# must stay text -- not be auto-coerced into numbers -- so force the cells
# to a text format before writing, then restore the default "Normal" style
# once the values are in place.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "43.603.77"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "2.322.52"
$ws.Range("E3").Value = "  +4.30%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "268.69"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").Value = "93.76"
$ws.Range("E6").Value = "  +8.18%  "
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +2.38%  "
$ws.Range("D10").Value = "44.53"
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("D11").Value = "0.0941"
$ws.Range("E11").Value = "  +2.39%  "
$ws.Range("D12").Value = "8.05"
$ws.Range("E12").Value = "  +4.80%  "
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "2.652.49"
$ws.Range("E14").Value = "  +3.14%  "
$ws.Range("D15").Value = "15.44"
$ws.Range("E15").Value = "  +4.11%  "
$ws.Range("D16").Value = "0.858"
$ws.Range("E16").Value = "  +8.80%  "
$ws.Range("D17").Value = "2.315.42"
$ws.Range("E17").Value = "  +3.95%  "
$ws.Range("D18").Value = "43.564.60"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("E19").Value = "  +4.47%  "
$ws.Range("E20").Value = "  +7.03%  "
$ws.Range("D21").Value = "71.41"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("D22").Value = "236.90"
$ws.Range("E22").Value = "  +1.92%  "
$ws.Range("E23").Value = "  -3.53%  "
$ws.Range("E24").Value = "  +10.73%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "11.30"
$ws.Range("E26").Value = "  +5.23%  "
$ws.Range("E27").Value = "  -1.06%  "
$ws.Range("D28").Value = "3.48"
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("D29").Value = "2.29"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").Value = "38.44"
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("D31").Value = "22.66"
$ws.Range("E31").Value = "  +9.56%  "
$ws.Range("D32").Value = "171.71"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("D34").Value = "5.47"
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("E35").Value = "  +2.52%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "0.0355"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.108"
$ws.Range("E37").Value = "  -2.04%  "
$ws.Range("D38").Value = "4.34"
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("D39").Value = "3.40"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("D40").Value = "2.36"
$ws.Range("E40").Value = "  +7.73%  "
$ws.Range("E41").Value = "  +14.34%  "
$ws.Range("D42").Value = "1.36"
$ws.Range("E42").Value = "  +20.25%  "
$ws.Range("D43").Value = "12.04"
$ws.Range("E43").Value = "  -2.63%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "9.06"
$ws.Range("E44").Value = "  +6.87%  "
$ws.Range("B45").Value = "MultiversX"
$ws.Range("C45").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D45").Value = "61.51"
$ws.Range("E45").Value = "  -2.95%  "
$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D46").Value = "5.37"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("E47").Value = "  +3.58%  "
$ws.Range("D48").Value = "100.19"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("E49").Value = "  +1.12%  "
$ws.Range("D50").Value = "2.533.37"
$ws.Range("E50").Value = "  +3.15%  "
$ws.Range("E51").Value = "  -0.17%  "

$priceRange.Style = "Normal"
